# Loan RBI, Variable Instalments
#
# Insert a new (blank) column into the "Repayment schedule" sheet just
# before the existing "Late" column, to make room for a new "Variable"
# instalments field, then leave the "Repayment schedule" tab active /
# selected, matching where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make this the active sheet (so the column-insert operates on it and the
# workbook remembers it as the last-viewed tab).
$ws.Activate() | Out-Null

# Insert a new blank column at N - shifts Late/heading/Outstanding (old
# N:P) one column to the right (new O:Q).
$ws.Columns("N").Insert() | Out-Null

# New column inherits the width of the column to its left ("In Advance").
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Leave the selection where the author ended up.
$ws.Range("R6").Select() | Out-Null
